$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 9484.087
$ws.Range("I33").Value = 13332.375
$ws.Range("K33").Value = 13332.375
$ws.Range("M33").Value = -13103.375

$ws.Range("H80").Value = 1828.2916
$ws.Range("I80").Value = 1354
$ws.Range("K80").Value = 4062
$ws.Range("M80").Value = -3064

$ws.Range("H83").Value = 1828.2916
$ws.Range("I83").Value = 1354
$ws.Range("K83").Value = 12186
$ws.Range("M83").Value = -7194

$ws.Range("H86").Value = 1542.7142
$ws.Range("I86").Value = 799.5
$ws.Range("K86").Value = 799.5
$ws.Range("M86").Value = 323.5

$ws.Range("H89").Value = 1542.7142
$ws.Range("I89").Value = 799.5
$ws.Range("K89").Value = 3997.5
$ws.Range("M89").Value = 1618.5

$ws.Range("M106").ClearContents()
$ws.Range("H106").Value = 2945
$ws.Range("I106").Value = 2945
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2945
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = -2314

$ws.Range("H107").Value = 1181
$ws.Range("I107").Value = 1060.2941
$ws.Range("J107").Value = 1474.1428
$ws.Range("K107").Value = 1060.2941
$ws.Range("L107").Value = 1474.1428
$ws.Range("M107").Value = 859.7058999999999
$ws.Range("N107").Value = -5314.1428

$ws.Range("H138").Value = 2171.125
$ws.Range("J138").Value = 3294.1155
$ws.Range("L138").Value = 9882.3465
$ws.Range("N138").Value = -20162.3465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 51655.633
$ws.Range("I132").Value = 2769.4888
$ws.Range("K132").Value = 8308.466400000001
$ws.Range("M132").Value = -5778.466400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5262.129
$ws.Range("I105").Value = 3154.1428
$ws.Range("J105").Value = 5876.9585
$ws.Range("K105").Value = 3154.1428
$ws.Range("L105").Value = 5876.9585
$ws.Range("M105").Value = -1407.1428
$ws.Range("N105").Value = -9370.958500000001

$ws.Range("H107").Value = 2039.7778
$ws.Range("I107").Value = 2008.2222
$ws.Range("J107").Value = 2102.889
$ws.Range("K107").Value = 2008.2222
$ws.Range("L107").Value = 2102.889
$ws.Range("M107").Value = -88.22219999999993
$ws.Range("N107").Value = -5942.889

$ws.Range("H132").Value = 118992.336
$ws.Range("J132").Value = 118992.336
$ws.Range("L132").Value = 118992.336
$ws.Range("N132").Value = -129112.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 209.26315
$ws.Range("I7").Value = 66
$ws.Range("K7").Value = 66
$ws.Range("M7").Value = 47

$ws.Range("H58").Value = 15280721
$ws.Range("I58").Value = 2941
$ws.Range("J58").Value = 68752950
$ws.Range("K58").Value = 2941
$ws.Range("L58").Value = 68752950
$ws.Range("M58").Value = -2738
$ws.Range("N58").Value = -68753356

$ws.Range("H99").Value = 5658.1665
$ws.Range("I99").Value = 5809.091
$ws.Range("K99").Value = 5809.091
$ws.Range("M99").Value = -4311.091

$ws.Range("H126").Value = 5658.1665
$ws.Range("I126").Value = 5809.091
$ws.Range("K126").Value = 17427.273
$ws.Range("M126").Value = -14957.273

$ws.Range("H136").Value = 15280721
$ws.Range("I136").Value = 2941
$ws.Range("J136").Value = 68752950
$ws.Range("K136").Value = 8823
$ws.Range("L136").Value = 206258850
$ws.Range("M136").Value = -6273
$ws.Range("N136").Value = -206263950

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 168.33333
$ws.Range("I2").Value = 198.9
$ws.Range("K2").Value = 1193.4
$ws.Range("M2").Value = -1080.4

$ws.Range("H15").Value = 350
$ws.Range("I15").Value = 50
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = -10
$ws.Range("N15").Value = -1780

$ws.Range("H68").Value = 7998.75
$ws.Range("I68").Value = 14997.5
$ws.Range("K68").Value = 44992.5
$ws.Range("M68").Value = -44181.5

$ws.Range("H71").Value = 7998.75
$ws.Range("I71").Value = 14997.5
$ws.Range("K71").Value = 134977.5
$ws.Range("M71").Value = -130921.5

$ws.Range("H113").Value = 1113.4615
$ws.Range("J113").Value = 2534
$ws.Range("L113").Value = 7602
$ws.Range("N113").Value = -11942

$ws.Range("H131").Value = 23810898
$ws.Range("J131").Value = 1847.3
$ws.Range("L131").Value = 5541.9
$ws.Range("N131").Value = -15621.9

$ws.Range("H139").Value = 2857.923
$ws.Range("I139").Value = 2644.2727
$ws.Range("K139").Value = 7932.8181
$ws.Range("M139").Value = -2792.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I35").Value = 18000
$ws.Range("K35").Value = 18000
$ws.Range("M35").Value = -17702

$ws.Range("H122").Value = 3368.9092
$ws.Range("I122").Value = 2294
$ws.Range("J122").Value = 5250
$ws.Range("K122").Value = 6882
$ws.Range("L122").Value = 15750
$ws.Range("M122").Value = -4432
$ws.Range("N122").Value = -20650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 946.52
$ws.Range("I16").Value = 1022.7
$ws.Range("J16").Value = 641.8
$ws.Range("K16").Value = 1022.7
$ws.Range("L16").Value = 641.8
$ws.Range("M16").Value = -852.7
$ws.Range("N16").Value = -981.8

$ws.Range("H40").Value = 7277.8
$ws.Range("I40").Value = 4200
$ws.Range("K40").Value = 4200
$ws.Range("M40").Value = -4064

$ws.Range("H122").Value = 4041.1052
$ws.Range("I122").Value = 3833.5356
$ws.Range("K122").Value = 11500.6068
$ws.Range("M122").Value = -9050.606800000001

$ws.Range("H132").Value = 2820.75
$ws.Range("I132").Value = 1614.2222
$ws.Range("K132").Value = 4842.6666
$ws.Range("M132").Value = -2312.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 6.25
$ws.Range("I12").Value = 6.181818
$ws.Range("K12").Value = 6.181818
$ws.Range("M12").Value = 135.818182

$ws.Range("H81").Value = 5109.59
$ws.Range("I81").Value = 3149.318
$ws.Range("J81").Value = 7646.4116
$ws.Range("K81").Value = 6298.636
$ws.Range("L81").Value = 15292.8232
$ws.Range("M81").Value = -5237.636
$ws.Range("N81").Value = -17414.8232

$ws.Range("H84").Value = 5109.59
$ws.Range("I84").Value = 3149.318
$ws.Range("J84").Value = 7646.4116
$ws.Range("K84").Value = 31493.18
$ws.Range("L84").Value = 76464.11600000001
$ws.Range("M84").Value = -26189.18
$ws.Range("N84").Value = -87072.11600000001

$ws.Range("H107").Value = 1617.6666
$ws.Range("I107").Value = 1207.8667
$ws.Range("K107").Value = 3623.6001
$ws.Range("M107").Value = -1703.6001

$ws.Range("H130").Value = 14999.75
$ws.Range("J130").Value = 14999.75
$ws.Range("L130").Value = 14999.75
$ws.Range("N130").Value = -25039.75

$ws.Range("H132").Value = 4411.1665
$ws.Range("I132").Value = 4268.45
$ws.Range("J132").Value = 5124.75
$ws.Range("K132").Value = 12805.35
$ws.Range("L132").Value = 15374.25
$ws.Range("M132").Value = -10275.35
$ws.Range("N132").Value = -20434.25

$ws.Range("H136").Value = 3847.5
$ws.Range("I136").Value = 2967.647
$ws.Range("J136").Value = 8833.333000000001
$ws.Range("K136").Value = 8902.940999999999
$ws.Range("L136").Value = 26499.999
$ws.Range("M136").Value = -6352.940999999999
$ws.Range("N136").Value = -31599.999
